# Generate Report for Handoff
# Updates the Overview/zh-cn/de-de sheets to reflect that the file is now
# "Ready for handoff" instead of "In Translation", and refreshes the
# corresponding handoff timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: status for both locales, and latest overall handoff date
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-03-24 04:37:09"

# zh-cn sheet: status + latest handoff datetime
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-24 04:37:05"

# de-de sheet: status + latest handoff datetime
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-24 04:37:09"
